$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every record.
# The whole column (row 2 through the last data row) was bumped by one day,
# from serial 45189 (2023-09-20) to 45190 (2023-09-21).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$ws.Range("C2:C$lastRow").Value = 45190
